
# ==========================================================================
# Edit script: "Asignacion de tareas para los demas colaboradores"
#
# 1. "Ing. Osorio " -> "Osorio "
# 2. "Ing. Nelson " -> "Nelson "
# 3. "Ing. Pedro de Jesús" -> "Pedro de Jesús" (drop the "Ing. " run,
#    keep "Pedro de " / "Jesús" as two separate runs)
# 4. New paragraphs "Diego" and "Marco " after the Pedro de Jesús line
# 5. Logotipo image paragraph gains bold/noProof run + paragraph-mark
#    formatting
# 6. Four new paragraphs after the logo image: "Rediseño del logotipo",
#    "Colaborador a cargo Marco  ", the HTTP/.com sentence (with
#    proofErr spell-check markers), "Colaborador a cargo Diego  "
# 7. "Fecha de termino:" date runs collapsed into a single run
# 8. Trim the 6 trailing empty paragraphs down to 1
# ==========================================================================

$d = $word.ActiveDocument

function Get-ParaByText($text, $exact = $true) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($exact) {
            if ($t -eq $text) { return $p }
        } else {
            if ($t -like "*$text*") { return $p }
        }
    }
    return $null
}

# --- Step 1: "Ing. Osorio " -> "Osorio " -----------------------------------
$pOsorio = Get-ParaByText "Ing. Osorio "
$pOsorio.Range.Find.Execute("Ing. Osorio ", $true, $false, $false, $false, $false, $true, 1, $false, "Osorio ", 2) | Out-Null

# --- Step 2: "Ing. Nelson " -> "Nelson " ------------------------------------
$pNelson = Get-ParaByText "Ing. Nelson "
$pNelson.Range.Find.Execute("Ing. Nelson ", $true, $false, $false, $false, $false, $true, 1, $false, "Nelson ", 2) | Out-Null

# --- Step 3: "Ing. Pedro de Jesús" -> "Pedro de Jesús" (2 runs preserved) ---
$pPedro = Get-ParaByText "Ing. Pedro de Jesús"
$xmlPedro = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:r><w:t xml:space="preserve">Pedro de </w:t></w:r><w:r><w:t>Jesús</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pPedro.Range.InsertXML($xmlPedro)

# --- Step 4: insert "Diego" and "Marco " paragraphs after it ---------------
$pPedro = Get-ParaByText "Pedro de Jesús"
$pPedro.Range.InsertParagraphAfter()
$pDiegoHolder = $pPedro.Next()
$xmlDiego = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:r><w:t>Diego</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pDiegoHolder.Range.InsertXML($xmlDiego)

$pDiego = Get-ParaByText "Diego"
$pDiego.Range.InsertParagraphAfter()
$pMarcoHolder = $pDiego.Next()
$xmlMarco = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:r><w:t xml:space="preserve">Marco </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pMarcoHolder.Range.InsertXML($xmlMarco)

# --- Step 5: bold + noProof on the logo image paragraph --------------------
$pImage = Get-ParaByText "Logotipo:"
$pImage = $pImage.Next()
$xmlImage = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="24CD3D28" wp14:editId="0D8B7DFE"><wp:extent cx="1221105" cy="1060450"/><wp:effectExtent l="0" t="0" r="0" b="6350"/><wp:docPr id="179283639" name="Imagen 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="179283639" name=""/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill rotWithShape="1"><a:blip r:embed="rId4"/><a:srcRect l="1357" t="1575"/><a:stretch/></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="1238633" cy="1075672"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:ln><a:noFill/></a:ln><a:extLst><a:ext uri="{53640926-AAD7-44D8-BBD7-CCE9431645EC}"><a14:shadowObscured xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"/></a:ext></a:extLst></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pImage.Range.InsertXML($xmlImage)

# --- Step 6: four new paragraphs after the logo image -----------------------
$pImage.Range.InsertParagraphAfter()
$pRedisenoHolder = $pImage.Next()
$xmlRediseno = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Rediseño del logotipo</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pRedisenoHolder.Range.InsertXML($xmlRediseno)

$pRediseno = Get-ParaByText "Rediseño del logotipo"
$pRediseno.Range.InsertParagraphAfter()
$pColabMarcoHolder = $pRediseno.Next()
$xmlColabMarco = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:r><w:t xml:space="preserve">Colaborador a cargo Marco  </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pColabMarcoHolder.Range.InsertXML($xmlColabMarco)

$pColabMarco = Get-ParaByText "Colaborador a cargo Marco  "
$pColabMarco.Range.InsertParagraphAfter()
$pSubidaHolder = $pColabMarco.Next()
$xmlSubida = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>La subida al servidor será mediante un servidor HTTP con terminación .</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>com</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pSubidaHolder.Range.InsertXML($xmlSubida)

$pSubida = Get-ParaByText "servidor HTTP" $false
$pSubida.Range.InsertParagraphAfter()
$pColabDiegoHolder = $pSubida.Next()
$xmlColabDiego = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:r><w:t xml:space="preserve">Colaborador a cargo Diego  </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pColabDiegoHolder.Range.InsertXML($xmlColabDiego)

# --- Step 7: collapse "Fecha de termino:" date runs into one ----------------
$pTermino = Get-ParaByText "termino" $false
$xmlTermino = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Fecha de termino:</w:t></w:r><w:r><w:t xml:space="preserve"> 14/08/2024 </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pTermino.Range.InsertXML($xmlTermino)

# --- Step 8: trim the 6 trailing empty paragraphs down to 1 -----------------
$pTermino = Get-ParaByText "termino" $false
$pFirstBlank = $pTermino.Next()
$total = $d.Paragraphs.Count
$pLastBlankButOne = $d.Paragraphs.Item($total - 1)
if ($pLastBlankButOne.Range.Start -gt $pFirstBlank.Range.Start) {
    $delRange = $d.Range($pFirstBlank.Range.Start, $pLastBlankButOne.Range.End)
    $delRange.Delete()
}

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host "$i : [$($p.Range.Text)]"
}
